$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.375.06"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "1.845.90"
$ws.Range("E3").Value = "  -0.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.40"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6295"
$ws.Range("E6").Value = "  +0.51%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -0.69%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.50"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07724"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("D12").Value = "1.871.84"
$ws.Range("E12").Value = "  -4.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.999"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6837"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000009995"
$ws.Range("E15").Value = "  +2.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "82.84"
$ws.Range("E16").Value = "  -1.07%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.149"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").Value = "29.424.35"
$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "227.65"
$ws.Range("E19").Value = "  -2.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.44"
$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("E21").Value = "  +0.04%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.549"
$ws.Range("E22").Value = "  -0.60%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.29"
$ws.Range("E24").Value = "  +1.25%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1398"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.362"
$ws.Range("E26").Value = "  -0.85%  "

$ws.Range("E27").Value = "  -0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.465"
$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05693"
$ws.Range("E29").Value = "  -3.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.257"
$ws.Range("E30").Value = "  -0.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.129"
$ws.Range("E31").Value = "  +0.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.015"
$ws.Range("E32").Value = "  -0.52%  "

$ws.Range("E33").Value = "  -2.32%  "

$ws.Range("E34").Value = "  -1.39%  "

$ws.Range("E35").Value = "  -0.67%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.589"
$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("D37").Value = "1.254.13"
$ws.Range("E37").Value = "  +1.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01813"
$ws.Range("E38").Value = "  +2.04%  "

$ws.Range("E39").Value = "  -0.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9130"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.211"
$ws.Range("E41").Value = "  +1.24%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").Value = "2.006.52"
$ws.Range("E43").Value = "  -4.73%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.10"
$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "66.37"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.046"
$ws.Range("E46").Value = "  -3.75%  "

$ws.Range("E47").Value = "  -0.16%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.116"
$ws.Range("E48").Value = "  -0.17%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.694"
$ws.Range("E49").Value = "  -1.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1127"
$ws.Range("E50").Value = "  +1.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05734"
$ws.Range("E51").Value = "  -0.33%  "
